# Natmi following Dr Hou advice
# Update LR-pair sheet (Efna5-Ephb1) with additional sending/target cluster "ECs"
# and recomputed values across all combinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order (A..T) matches headers already present in row 1.

$data = @(
    @("ECs",  "Efna5", "Ephb1", "ECs", 2, 0.6666666666666666, 0.230855, 0.692565, 0.06377305075821572, 0.06377305075821572, 2, 0.6666666666666666, 1.471191666666667, 4.413575, 0.6447353255635294, 0.6447353255635294, 0.3396319522083334, 3.056687569875, 0.0411167386427777, 0.0411167386427777),
    @("ECs",  "Efna5", "Ephb1", "sCs", 2, 0.6666666666666666, 0.230855, 0.692565, 0.06377305075821572, 0.06377305075821572, 3, 1, 0.8106620000000001, 2.431986, 0.3552646744364706, 0.3552646744364706, 0.18714537601, 1.68430838409, 0.02265631211543803, 0.02265631211543803),
    @("FAPs", "Efna5", "Ephb1", "ECs", 3, 1, 2.900731333333333, 8.702194, 0.8013189515350044, 0.8013189515350045, 2, 0.6666666666666666, 1.471191666666667, 4.413575, 0.6447353255635294, 0.6447353255635294, 4.267531764838889, 38.40778588355, 0.5166386350981471, 0.5166386350981472),
    @("FAPs", "Efna5", "Ephb1", "sCs", 3, 1, 2.900731333333333, 8.702194, 0.8013189515350044, 0.8013189515350045, 3, 1, 0.8106620000000001, 2.431986, 0.3552646744364706, 0.3552646744364706, 2.351512664142667, 21.163613977284, 0.2846803164368573, 0.2846803164368574),
    @("sCs",  "Efna5", "Ephb1", "ECs", 3, 1, 0.4883596666666667, 1.465079, 0.1349079977067798, 0.1349079977067798, 2, 0.6666666666666666, 1.471191666666667, 4.413575, 0.6447353255635294, 0.6447353255635294, 0.7184706719361111, 6.466236047425, 0.08697995182260455, 0.08697995182260455),
    @("sCs",  "Efna5", "Ephb1", "sCs", 3, 1, 0.4883596666666667, 1.465079, 0.1349079977067798, 0.1349079977067798, 3, 1, 0.8106620000000001, 2.431986, 0.3552646744364706, 0.3552646744364706, 0.3958946240993334, 3.563051616894001, 0.04792804588417526, 0.04792804588417526)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
